# Refresh the crypto price/volume snapshot for this run.
# (Coin/Link/Price/Volume columns are B/C/D/E; column A is a
#  static row index and is left untouched.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.613.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.493.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.43%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.64%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.492.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.32%  "

$ws.Range("E9").Value = "  +2.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.34%  "

$ws.Range("E11").Value = "  +4.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.436"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.097.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.35%  "

$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.46%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.48%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.537.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.483.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.62%  "

$ws.Range("E19").Value = "  +3.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.65%  "

$ws.Range("E22").Value = "  +2.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000124"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.76%  "

$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.528"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.02%  "

$ws.Range("E28").Value = "  +1.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.03%  "

$ws.Range("E30").Value = "  +5.22%  "

$ws.Range("E31").Value = "  +6.44%  "

$ws.Range("E32").Value = "  +3.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.32%  "

$ws.Range("E34").Value = "  +5.93%  "

$ws.Range("E35").Value = "  +0.13%  "

$ws.Range("E36").Value = "  +1.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.903"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.58%  "

$ws.Range("E39").Value = "  +5.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0744"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.42"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.811.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.95%  "

$ws.Range("E48").Value = "  +3.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "356.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.39%  "
